# Update the "path" column (A) of the index sheet so each sample file name
# becomes the full absolute path used on the author's machine, and reorder
# the hic/rnaseq groupings to match. Also resize column A to fit the much
# longer paths and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "/Users/scottronquist/projects/4DNvestigator/data/projects/myod/raw/hic/Sample_64585_trim.hic"
$ws.Range("A3").Value = "/Users/scottronquist/projects/4DNvestigator/data/projects/myod/raw/hic/Sample_71530_trim.hic"
$ws.Range("A4").Value = "/Users/scottronquist/projects/4DNvestigator/data/projects/myod/raw/hic/Sample_71531_trim.hic"
$ws.Range("A5").Value = "/Users/scottronquist/projects/4DNvestigator/data/projects/myod/raw/rnaseq/Sample_63246_rsem.genes.results"
$ws.Range("A6").Value = "/Users/scottronquist/projects/4DNvestigator/data/projects/myod/raw/rnaseq/Sample_63247_rsem.genes.results"
$ws.Range("A7").Value = "/Users/scottronquist/projects/4DNvestigator/data/projects/myod/raw/rnaseq/Sample_63248_rsem.genes.results"
$ws.Range("A8").Value = "/Users/scottronquist/projects/4DNvestigator/data/projects/myod/raw/rnaseq/Sample_63249_rsem.genes.results"
$ws.Range("A9").Value = "/Users/scottronquist/projects/4DNvestigator/data/projects/myod/raw/rnaseq/Sample_63250_rsem.genes.results"
$ws.Range("A10").Value = "/Users/scottronquist/projects/4DNvestigator/data/projects/myod/raw/rnaseq/Sample_63251_rsem.genes.results"
$ws.Range("A11").Value = "/Users/scottronquist/projects/4DNvestigator/data/projects/myod/raw/rnaseq/Sample_63273_rsem.genes.results"
$ws.Range("A12").Value = "/Users/scottronquist/projects/4DNvestigator/data/projects/myod/raw/rnaseq/Sample_63274_rsem.genes.results"
$ws.Range("A13").Value = "/Users/scottronquist/projects/4DNvestigator/data/projects/myod/raw/rnaseq/Sample_63275_rsem.genes.results"

$ws.Columns.Item(1).ColumnWidth = 152.66666666666666
$ws.Range("A5").Select() | Out-Null
